$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsAdmin = $wb.Worksheets.Item("AdminLogin")

# Create the new "UserData" sheet as a copy of "Login" (placed after "AdminLogin"),
# which preserves sheet formatting (row height, columns, fonts/styles, hyperlink
# formatting) identical to the existing Login sheet.
$wsLogin.Copy([System.Type]::Missing, $wsAdmin)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "UserData"

# Remove the hyperlink that was copied over from Login (it still points at
# raviuser@yopmail.com) so we can replace it with the new sub-user address.
$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = "autouser3@yopmail.com"
$ws3.Range("B2").Value = 12345678

$ws3.Hyperlinks.Add($ws3.Range("A2"), "mailto:autouser3@yopmail.com")

# Re-apply the Login A2 cell format so the hyperlink cell keeps using the
# existing shared "Hyperlink" style instead of a newly synthesized one.
$wsLogin.Range("A2").Copy()
$ws3.Range("A2").PasteSpecial(-4122)
$ws3.Range("A2").Value = "autouser3@yopmail.com"

# Match the column widths used on the new sheet (closest achievable given the
# 1/6-character granularity of the ColumnWidth COM property).
$ws3.Columns.Item(1).ColumnWidth = 45
$ws3.Columns.Item(2).ColumnWidth = 42.166666666666664

# New sheet opens with A2 selected/active.
$ws3.Range("A2").Select()

# The Login sheet is no longer the active tab; its lingering selection moves
# to A1:B2.
$wsLogin.Range("A1:B2").Select()

$ws3.Activate()
